# Applies the "Updated symbol list" data refresh described in the commit:
#   - rows 8-18: a new "GateToken" entry is inserted at row 8 (pulled up from
#     its old row 18 slot with refreshed price/volume), shifting the other
#     exchange-token rows (BTSEToken, MXToken, Liechtenstein..., WazirX, ...,
#     LEO) down by one position each
#   - every row keeps the same coin/link but gets refreshed Price (col D) and
#     Volume(1h) (col E) figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link columns (B, C): plain text, safe to assign directly ---
$ws.Range('B8').Value = 'GateToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('B10').Value = 'MXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('B12').Value = 'WazirX'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

# --- Price / Volume(1h) columns (D, E): these look numeric ("329.83",
# "2.12%", "--", "0.00000000751", ...) but must stay plain text, matching
# the source sheet (t="inlineStr" everywhere, no numeric cells in D/E).
# A leading apostrophe forces Excel to store the literal as text instead of
# coercing it into a number/percentage.
$ws.Range('D2').Value = "'329.83"
$ws.Range('E2').Value = "'2.12%"
$ws.Range('D3').Value = "'41.03"
$ws.Range('E3').Value = "'3.32%"
$ws.Range('D4').Value = "'5.643"
$ws.Range('E4').Value = "'-4.04%"
$ws.Range('D5').Value = "'0.08162"
$ws.Range('E5').Value = "'1.65%"
$ws.Range('D6').Value = "'2.039"
$ws.Range('E6').Value = "'5.29%"
$ws.Range('D7').Value = "'8.745"
$ws.Range('E7').Value = "'1.03%"
$ws.Range('D8').Value = "'4.543"
$ws.Range('E8').Value = "'-0.63%"
$ws.Range('D9').Value = "'2.953"
$ws.Range('E9').Value = "'0.11%"
$ws.Range('D10').Value = "'0.9179"
$ws.Range('E10').Value = "'-1.33%"
$ws.Range('D11').Value = "'0.1254"
$ws.Range('E11').Value = "'2.42%"
$ws.Range('D12').Value = "'0.1951"
$ws.Range('E12').Value = "'-0.55%"
$ws.Range('D13').Value = "'0.09367"
$ws.Range('E13').Value = "'2.81%"
$ws.Range('D14').Value = "'0.03710"
$ws.Range('E14').Value = "'4.93%"
$ws.Range('D15').Value = "'0.1054"
$ws.Range('E15').Value = "'10.21%"
$ws.Range('D16').Value = "'0.001298"
$ws.Range('E16').Value = "'0.03%"
$ws.Range('D17').Value = "'0.006145"
$ws.Range('E17').Value = "'0.32%"
$ws.Range('D18').Value = "'3.432"
$ws.Range('E18').Value = "'2.37%"
$ws.Range('E19').Value = "'-2.12%"
$ws.Range('D20').Value = "'8.256"
$ws.Range('E20').Value = "'-5.66%"
$ws.Range('E21').Value = "'-1.67%"
$ws.Range('D22').Value = "'0.2653"
$ws.Range('E22').Value = "'10.18%"
$ws.Range('D23').Value = "'0.04429"
$ws.Range('E23').Value = "'0.57%"
$ws.Range('D24').Value = "'0.001266"
$ws.Range('E24').Value = "'0.42%"
$ws.Range('D25').Value = "'0.004293"
$ws.Range('E25').Value = "'-2.05%"
$ws.Range('D26').Value = "'0.0001182"
$ws.Range('E26').Value = "'3.71%"
$ws.Range('D39').Value = "'0.02753"
$ws.Range('E39').Value = "'13.56%"
$ws.Range('D40').Value = "'0.05441"
$ws.Range('E40').Value = "'4.15%"
$ws.Range('D41').Value = "'0.007654"
$ws.Range('E41').Value = "'2.77%"
$ws.Range('D42').Value = "'0.009437"
$ws.Range('E42').Value = "'0.70%"
$ws.Range('E43').Value = "'0.77%"
$ws.Range('D44').Value = "'0.002114"
$ws.Range('E44').Value = "'-0.28%"
$ws.Range('D45').Value = "'0.01203"
$ws.Range('E45').Value = "'6.86%"
$ws.Range('D46').Value = "'0.00006883"
$ws.Range('E46').Value = "'2.42%"
$ws.Range('D47').Value = "'0.00000000751"
$ws.Range('E47').Value = "'0.20%"
$ws.Range('E48').Value = "'60.53%"
$ws.Range('D49').Value = "'0.003532"
$ws.Range('E49').Value = "'17.69%"
$ws.Range('D50').Value = "'0.00002104"
$ws.Range('E50').Value = "'0.20%"
$ws.Range('E51').Value = "'0.20%"

# Quote-prefixing a numeric-looking literal makes Excel silently stamp the
# cell with a Text number format (so it keeps round-tripping as text). The
# source cells carry no explicit style at all, so strip that incidental
# formatting back off -- this only touches number format/style, not the
# text values just written above.
$ws.Range('D2:E51').ClearFormats()
